# Update column C (Förändrad) for rows 2-9 from serial date 45184 (2023-09-15)
# to 45185 (2023-09-16), leaving everything else untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 9; $row++) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    if ($cell.Value2 -eq 45184) {
        $cell.Value2 = 45185
    }
}
